$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new log rows (84 and 85) to the feed logs sheet
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = 1
$ws.Range("C84").Value = "2024-06-16 16:14:04"
$ws.Range("D84").Value = 200
$ws.Range("E84").Value = 14

$ws.Range("A85").Value = 84
$ws.Range("B85").Value = 2
$ws.Range("C85").Value = "2024-06-16 16:14:04"
$ws.Range("D85").Value = 200
$ws.Range("E85").Value = 0
